# Update the LR-pairs data table (Efnb3-Epha4) with new TPM-derived values.
# The sending-cluster set now also includes "ECs", so the table grows from
# 8 data rows (rows 2-9) to 12 data rows (rows 2-13): 3 sending clusters
# (ECs, FAPs, MuSCs) x 4 target clusters (ECs, FAPs, MuSCs, Resolving-Mac).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efnb3"
$ws.Range("C2").Value = "Epha4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2151756666666667
$ws.Range("H2").Value = 0.645527
$ws.Range("I2").Value = 0.1791915537270423
$ws.Range("J2").Value = 0.1791915537270423
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 10.23061133333333
$ws.Range("N2").Value = 30.691834
$ws.Range("O2").Value = 0.4855635428718841
$ws.Range("P2").Value = 0.4855635428718841
$ws.Range("Q2").Value = 2.201378614057556
$ws.Range("R2").Value = 19.812407526518
$ws.Range("S2").Value = 0.08700888568042021
$ws.Range("T2").Value = 0.08700888568042021

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efnb3"
$ws.Range("C3").Value = "Epha4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2151756666666667
$ws.Range("H3").Value = 0.645527
$ws.Range("I3").Value = 0.1791915537270423
$ws.Range("J3").Value = 0.1791915537270423
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 8.775186333333332
$ws.Range("N3").Value = 26.325559
$ws.Range("O3").Value = 0.4164864079521221
$ws.Range("P3").Value = 0.4164864079521222
$ws.Range("Q3").Value = 1.888206569399222
$ws.Range("R3").Value = 16.993859124593
$ws.Range("S3").Value = 0.07463084654713555
$ws.Range("T3").Value = 0.07463084654713555

# Row 4: ECs -> MuSCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efnb3"
$ws.Range("C4").Value = "Epha4"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2151756666666667
$ws.Range("H4").Value = 0.645527
$ws.Range("I4").Value = 0.1791915537270423
$ws.Range("J4").Value = 0.1791915537270423
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.034752
$ws.Range("N4").Value = 6.104255999999999
$ws.Range("O4").Value = 0.09657305490303886
$ws.Range("P4").Value = 0.09657305490303887
$ws.Range("Q4").Value = 0.4378291181013332
$ws.Range("R4").Value = 3.940462062911999
$ws.Range("S4").Value = 0.01730507575624249
$ws.Range("T4").Value = 0.0173050757562425

# Row 5: ECs -> Resolving-Mac
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Efnb3"
$ws.Range("C5").Value = "Epha4"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2151756666666667
$ws.Range("H5").Value = 0.645527
$ws.Range("I5").Value = 0.1791915537270423
$ws.Range("J5").Value = 0.1791915537270423
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.02901266666666667
$ws.Range("N5").Value = 0.087038
$ws.Range("O5").Value = 0.001376994272954919
$ws.Range("P5").Value = 0.001376994272954919
$ws.Range("Q5").Value = 0.006242819891777778
$ws.Range("R5").Value = 0.056185379026
$ws.Range("S5").Value = 0.0002467457432440308
$ws.Range("T5").Value = 0.0002467457432440308

# Row 6: FAPs -> ECs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efnb3"
$ws.Range("C6").Value = "Epha4"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.5010686666666667
$ws.Range("H6").Value = 1.503206
$ws.Range("I6").Value = 0.4172742870736815
$ws.Range("J6").Value = 0.4172742870736815
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 10.23061133333333
$ws.Range("N6").Value = 30.691834
$ws.Range("O6").Value = 0.4855635428718841
$ws.Range("P6").Value = 0.4855635428718841
$ws.Range("Q6").Value = 5.126238779978223
$ws.Range("R6").Value = 46.136149019804
$ws.Range("S6").Value = 0.2026131811808364
$ws.Range("T6").Value = 0.2026131811808364

# Row 7: FAPs -> FAPs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efnb3"
$ws.Range("C7").Value = "Epha4"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.5010686666666667
$ws.Range("H7").Value = 1.503206
$ws.Range("I7").Value = 0.4172742870736815
$ws.Range("J7").Value = 0.4172742870736815
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 8.775186333333332
$ws.Range("N7").Value = 26.325559
$ws.Range("O7").Value = 0.4164864079521221
$ws.Range("P7").Value = 0.4164864079521222
$ws.Range("Q7").Value = 4.396970915794888
$ws.Range("R7").Value = 39.572738242154
$ws.Range("S7").Value = 0.1737890689541002
$ws.Range("T7").Value = 0.1737890689541002

# Row 8: FAPs -> MuSCs
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Efnb3"
$ws.Range("C8").Value = "Epha4"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.5010686666666667
$ws.Range("H8").Value = 1.503206
$ws.Range("I8").Value = 0.4172742870736815
$ws.Range("J8").Value = 0.4172742870736815
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.034752
$ws.Range("N8").Value = 6.104255999999999
$ws.Range("O8").Value = 0.09657305490303886
$ws.Range("P8").Value = 0.09657305490303887
$ws.Range("Q8").Value = 1.019550471637333
$ws.Range("R8").Value = 9.175954244735999
$ws.Range("S8").Value = 0.04029745263519304
$ws.Range("T8").Value = 0.04029745263519304

# Row 9: FAPs -> Resolving-Mac
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Efnb3"
$ws.Range("C9").Value = "Epha4"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.5010686666666667
$ws.Range("H9").Value = 1.503206
$ws.Range("I9").Value = 0.4172742870736815
$ws.Range("J9").Value = 0.4172742870736815
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.02901266666666667
$ws.Range("N9").Value = 0.087038
$ws.Range("O9").Value = 0.001376994272954919
$ws.Range("P9").Value = 0.001376994272954919
$ws.Range("Q9").Value = 0.01453733820311111
$ws.Range("R9").Value = 0.130836043828
$ws.Range("S9").Value = 0.0005745843035518059
$ws.Range("T9").Value = 0.0005745843035518059

# Row 10: MuSCs -> ECs
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Efnb3"
$ws.Range("C10").Value = "Epha4"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.4845693333333334
$ws.Range("H10").Value = 1.453708
$ws.Range("I10").Value = 0.4035341591992763
$ws.Range("J10").Value = 0.4035341591992763
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 10.23061133333333
$ws.Range("N10").Value = 30.691834
$ws.Range("O10").Value = 0.4855635428718841
$ws.Range("P10").Value = 0.4855635428718841
$ws.Range("Q10").Value = 4.957440513385778
$ws.Range("R10").Value = 44.616964620472
$ws.Range("S10").Value = 0.1959414760106275
$ws.Range("T10").Value = 0.1959414760106275

# Row 11: MuSCs -> FAPs
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Efnb3"
$ws.Range("C11").Value = "Epha4"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.4845693333333334
$ws.Range("H11").Value = 1.453708
$ws.Range("I11").Value = 0.4035341591992763
$ws.Range("J11").Value = 0.4035341591992763
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 8.775186333333332
$ws.Range("N11").Value = 26.325559
$ws.Range("O11").Value = 0.4164864079521221
$ws.Range("P11").Value = 0.4164864079521222
$ws.Range("Q11").Value = 4.252186191419111
$ws.Range("R11").Value = 38.269675722772
$ws.Range("S11").Value = 0.1680664924508864
$ws.Range("T11").Value = 0.1680664924508864

# Row 12: MuSCs -> MuSCs
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Efnb3"
$ws.Range("C12").Value = "Epha4"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.4845693333333334
$ws.Range("H12").Value = 1.453708
$ws.Range("I12").Value = 0.4035341591992763
$ws.Range("J12").Value = 0.4035341591992763
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.034752
$ws.Range("N12").Value = 6.104255999999999
$ws.Range("O12").Value = 0.09657305490303886
$ws.Range("P12").Value = 0.09657305490303887
$ws.Range("Q12").Value = 0.9859784201386665
$ws.Range("R12").Value = 8.873805781247999
$ws.Range("S12").Value = 0.03897052651160333
$ws.Range("T12").Value = 0.03897052651160333

# Row 13: MuSCs -> Resolving-Mac
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Efnb3"
$ws.Range("C13").Value = "Epha4"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.4845693333333334
$ws.Range("H13").Value = 1.453708
$ws.Range("I13").Value = 0.4035341591992763
$ws.Range("J13").Value = 0.4035341591992763
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.02901266666666667
$ws.Range("N13").Value = 0.087038
$ws.Range("O13").Value = 0.001376994272954919
$ws.Range("P13").Value = 0.001376994272954919
$ws.Range("Q13").Value = 0.01405864854488889
$ws.Range("R13").Value = 0.126527836904
$ws.Range("S13").Value = 0.0005556642261590818
$ws.Range("T13").Value = 0.0005556642261590818

